$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily TGP (terminal gate pricing) rollover:
# - old "today" (effective date 46002) rows become the new "yesterday" rows (46002 -> stays displayed as prior day)
# - new "today" rows get new effective date 46003 and freshly updated prices

# Row 8
$ws.Range("A8").Value = 46003
$ws.Range("D8").Value = 160.38999999999999
$ws.Range("E8").Value = 158.88
$ws.Range("F8").Value = 168.88
$ws.Range("G8").Value = 159

# Row 9
$ws.Range("A9").Value = 46003
$ws.Range("D9").Value = 160.38999999999999
$ws.Range("E9").Value = 158.88
$ws.Range("F9").Value = 168.88
$ws.Range("G9").Value = 159

# Row 10
$ws.Range("A10").Value = 46003
$ws.Range("D10").Value = 162.18
$ws.Range("E10").Value = 161.29
$ws.Range("F10").Value = 171.29
$ws.Range("G10").Value = 161.81

# Row 11
$ws.Range("A11").Value = 46002
$ws.Range("D11").Value = 160.96
$ws.Range("E11").Value = 159.43
$ws.Range("F11").Value = 169.43
$ws.Range("G11").Value = 159.54

# Row 12
$ws.Range("A12").Value = 46002
$ws.Range("D12").Value = 160.96
$ws.Range("E12").Value = 159.43
$ws.Range("F12").Value = 169.43
$ws.Range("G12").Value = 159.54

# Row 13
$ws.Range("A13").Value = 46002
$ws.Range("D13").Value = 162.79
$ws.Range("E13").Value = 161.66
$ws.Range("F13").Value = 171.66
$ws.Range("G13").Value = 162.18

# Row 17
$ws.Range("A17").Value = 46003
$ws.Range("D17").Value = 165.12
$ws.Range("E17").Value = 164.78
$ws.Range("F17").Value = 174.78

# Row 18
$ws.Range("A18").Value = 46002
$ws.Range("D18").Value = 165.73
$ws.Range("E18").Value = 165.36
$ws.Range("F18").Value = 175.36

# Row 22
$ws.Range("A22").Value = 46003
$ws.Range("D22").Value = 161.09
$ws.Range("E22").Value = 160.5
$ws.Range("F22").Value = 170.1
$ws.Range("G22").Value = 161.66

# Row 23
$ws.Range("A23").Value = 46003
$ws.Range("D23").Value = 167.2
$ws.Range("E23").Value = 165.4
$ws.Range("F23").Value = 175.4

# Row 24
$ws.Range("A24").Value = 46003
$ws.Range("D24").Value = 166.97
$ws.Range("E24").Value = 165.72
$ws.Range("F24").Value = 175.72

# Row 25
$ws.Range("A25").Value = 46003
$ws.Range("D25").Value = 167.58
$ws.Range("E25").Value = 165.15
$ws.Range("F25").Value = 175.15
$ws.Range("G25").Value = 164.92

# Row 26
$ws.Range("A26").Value = 46003
$ws.Range("D26").Value = 166.38
$ws.Range("E26").Value = 166.59
$ws.Range("F26").Value = 176.59

# Row 27
$ws.Range("A27").Value = 46002
$ws.Range("D27").Value = 162.1
$ws.Range("E27").Value = 160.97
$ws.Range("F27").Value = 170.57
$ws.Range("G27").Value = 162.13

# Row 28
$ws.Range("A28").Value = 46002
$ws.Range("D28").Value = 167.82
$ws.Range("E28").Value = 165.88
$ws.Range("F28").Value = 175.88

# Row 29
$ws.Range("A29").Value = 46002
$ws.Range("D29").Value = 167.59
$ws.Range("E29").Value = 166.2
$ws.Range("F29").Value = 176.2

# Row 30
$ws.Range("A30").Value = 46002
$ws.Range("D30").Value = 168.2
$ws.Range("E30").Value = 165.63
$ws.Range("F30").Value = 175.63
$ws.Range("G30").Value = 165.4

# Row 31
$ws.Range("A31").Value = 46002
$ws.Range("D31").Value = 166.99
$ws.Range("E31").Value = 167.07
$ws.Range("F31").Value = 177.07

# Row 35
$ws.Range("A35").Value = 46003
$ws.Range("D35").Value = 160.57
$ws.Range("E35").Value = 158.19999999999999
$ws.Range("F35").Value = 167.2

# Row 36
$ws.Range("A36").Value = 46002
$ws.Range("D36").Value = 161.08000000000001
$ws.Range("E36").Value = 158.68
$ws.Range("F36").Value = 167.68

# Row 40
$ws.Range("A40").Value = 46003
$ws.Range("D40").Value = 166.46
$ws.Range("E40").Value = 165.43
$ws.Range("F40").Value = 175.43

# Row 41
$ws.Range("A41").Value = 46003
$ws.Range("D41").Value = 166.17
$ws.Range("E41").Value = 165.85
$ws.Range("F41").Value = 175.85

# Row 42
$ws.Range("A42").Value = 46002
$ws.Range("D42").Value = 167.08
$ws.Range("E42").Value = 166.03
$ws.Range("F42").Value = 176.03

# Row 43
$ws.Range("A43").Value = 46002
$ws.Range("D43").Value = 166.8
$ws.Range("E43").Value = 166.45
$ws.Range("F43").Value = 176.45

# Row 47
$ws.Range("A47").Value = 46003
$ws.Range("D47").Value = 161.61000000000001
$ws.Range("E47").Value = 160.44
$ws.Range("F47").Value = 170.44

# Row 48
$ws.Range("A48").Value = 46003
$ws.Range("D48").Value = 161.41999999999999
$ws.Range("E48").Value = 160.53
$ws.Range("F48").Value = 170.53

# Row 49
$ws.Range("A49").Value = 46002
$ws.Range("D49").Value = 162.52000000000001
$ws.Range("E49").Value = 161.30000000000001
$ws.Range("F49").Value = 171.3

# Row 50
$ws.Range("A50").Value = 46002
$ws.Range("D50").Value = 162.32
$ws.Range("E50").Value = 161.38999999999999
$ws.Range("F50").Value = 171.39

# Row 54
$ws.Range("A54").Value = 46003
$ws.Range("D54").Value = 175.96
$ws.Range("E54").Value = 175.95
$ws.Range("F54").Value = 185.96

# Row 55
$ws.Range("A55").Value = 46003
$ws.Range("D55").Value = 164.14
$ws.Range("E55").Value = 171.07
$ws.Range("F55").Value = 181.07

# Row 56
$ws.Range("A56").Value = 46003
$ws.Range("D56").Value = 166.73

# Row 57
$ws.Range("A57").Value = 46003
$ws.Range("D57").Value = 165.74
$ws.Range("E57").Value = 165.34

# Row 58
$ws.Range("A58").Value = 46003
$ws.Range("D58").Value = 161.63999999999999
$ws.Range("E58").Value = 161.38999999999999
$ws.Range("F58").Value = 171.39

# Row 59
$ws.Range("A59").Value = 46003
$ws.Range("D59").Value = 168.03
$ws.Range("E59").Value = 173.21

# Row 60
$ws.Range("A60").Value = 46002
$ws.Range("D60").Value = 176.59
$ws.Range("E60").Value = 176.46
$ws.Range("F60").Value = 186.46

# Row 61
$ws.Range("A61").Value = 46002
$ws.Range("D61").Value = 164.76
$ws.Range("E61").Value = 171.54
$ws.Range("F61").Value = 181.54

# Row 62
$ws.Range("A62").Value = 46002
$ws.Range("D62").Value = 167.35

# Row 63
$ws.Range("A63").Value = 46002
$ws.Range("D63").Value = 166.34
$ws.Range("E63").Value = 165.81

# Row 64
$ws.Range("A64").Value = 46002
$ws.Range("D64").Value = 162.24
$ws.Range("E64").Value = 161.86000000000001
$ws.Range("F64").Value = 171.86

# Row 65
$ws.Range("A65").Value = 46002
$ws.Range("D65").Value = 168.64
$ws.Range("E65").Value = 173.71
